$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Unprotect()

$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.8440386337004069
$ws.Range("E2").Value = 0.005691372396843697

$ws.Range("D3").Value = 0.1559613662995931
$ws.Range("E3").Value = 0.008213552361396204

$ws.Range("E4").Value = 0.006084735030168931

$ws.Protect()
